$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Tipo"
$ws.Range("H2").Value = "Secundario"
$ws.Range("H21").Value = "Vencido"
$ws.Range("H6").Value = "Primario"

$ws.Range("H3").Value = "Secundario"
$ws.Range("H4").Value = "Secundario"
$ws.Range("H5").Value = "Secundario"
$ws.Range("H7").Value = "Secundario"
$ws.Range("H8").Value = "Secundario"
$ws.Range("H9").Value = "Primario"
$ws.Range("H10").Value = "Secundario"
$ws.Range("H11").Value = "Secundario"
$ws.Range("H12").Value = "Secundario"
$ws.Range("H13").Value = "Secundario"
$ws.Range("H14").Value = "Secundario"
$ws.Range("H15").Value = "Secundario"
$ws.Range("H16").Value = "Secundario"
$ws.Range("H17").Value = "Secundario"
$ws.Range("H18").Value = "Secundario"
$ws.Range("H19").Value = "Secundario"
$ws.Range("H20").Value = "Secundario"
$ws.Range("H22").Value = "Primario"
